$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 292, shifting existing rows 292..362 down to 293..363
$ws.Rows.Item(292).Insert()

# Populate the new row 292 with the new price observation
$ws.Cells.Item(292, 1).Value = 11
$ws.Cells.Item(292, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(292, 3).Value = "Bíobío"
$ws.Cells.Item(292, 4).Value = 44855
$ws.Cells.Item(292, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(292, 5).Value = 8
$ws.Cells.Item(292, 6).Value = "Fruta"
$ws.Cells.Item(292, 7).Value = 100102
$ws.Cells.Item(292, 8).Value = "Cítricos"
$ws.Cells.Item(292, 9).Value = 100102005
$ws.Cells.Item(292, 10).Value = "Naranja"
$ws.Cells.Item(292, 11).Value = "Lane Late"
$ws.Cells.Item(292, 12).Value = "Primera"
$ws.Cells.Item(292, 13).Value = 250
$ws.Cells.Item(292, 14).Value = 8000
$ws.Cells.Item(292, 15).Value = 8500
$ws.Cells.Item(292, 16).Value = 8240
$ws.Cells.Item(292, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(292, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(292, 19).Value = 549
$ws.Cells.Item(292, 20).Value = 15
